$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "92.354.90"
$ws.Range("E2").Value = "  -5.46%  "

$ws.Range("D3").Value = "3.317.70"
$ws.Range("E3").Value = "  -4.95%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'228.59"
$ws.Range("E5").Value = "  -8.24%  "

$ws.Range("D6").Value = "'615.01"
$ws.Range("E6").Value = "  -6.61%  "

$ws.Range("E7").Value = "  -5.80%  "

$ws.Range("D8").Value = "'0.377"
$ws.Range("E8").Value = "  -9.84%  "

$ws.Range("E9").Value = "  +0.09%  "

$ws.Range("D10").Value = "'0.914"
$ws.Range("E10").Value = "  -10.07%  "

$ws.Range("D11").Value = "3.318.95"
$ws.Range("E11").Value = "  -4.88%  "

$ws.Range("D12").Value = "'41.47"
$ws.Range("E12").Value = "  -6.62%  "

$ws.Range("D14").Value = "'5.92"
$ws.Range("E14").Value = "  -4.26%  "

$ws.Range("D15").Value = "91.925.68"
$ws.Range("E15").Value = "  -5.78%  "

$ws.Range("D16").Value = "3.934.19"
$ws.Range("E16").Value = "  -5.07%  "

$ws.Range("E17").Value = "  -6.63%  "

$ws.Range("E18").Value = "  -9.77%  "

$ws.Range("D19").Value = "3.316.73"
$ws.Range("E19").Value = "  -4.86%  "

$ws.Range("D20").Value = "'16.81"
$ws.Range("E20").Value = "  -9.31%  "

$ws.Range("D21").Value = "'10.86"
$ws.Range("E21").Value = "  -9.92%  "

$ws.Range("D22").Value = "'483.22"
$ws.Range("E22").Value = "  -6.94%  "

$ws.Range("D23").Value = "'3.22"
$ws.Range("E23").Value = "  -3.71%  "

$ws.Range("D24").Value = "'0.442"
$ws.Range("E24").Value = "  -11.46%  "

$ws.Range("E25").Value = "  -9.43%  "

$ws.Range("D26").Value = "'6.12"
$ws.Range("E26").Value = "  -9.62%  "

$ws.Range("D27").Value = "'89.20"
$ws.Range("E27").Value = "  -7.81%  "

$ws.Range("D28").Value = "'11.45"
$ws.Range("E28").Value = "  -8.67%  "

$ws.Range("E29").Value = "  -0.28%  "

$ws.Range("D30").Value = "'11.05"
$ws.Range("E30").Value = "  -10.43%  "

$ws.Range("E31").Value = "  -5.32%  "

$ws.Range("E32").Value = "  -8.91%  "

$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "  +0.06%  "

$ws.Range("D34").Value = "'0.169"
$ws.Range("E34").Value = "  -9.69%  "

$ws.Range("D35").Value = "'28.01"
$ws.Range("E35").Value = "  -9.76%  "

$ws.Range("E36").Value = "  -12.81%  "

$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D37").Value = "'520.00"
$ws.Range("E37").Value = "  +0.27%  "

$ws.Range("B38").Value = "USDe"
$ws.Range("C38").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = "  -0.02%  "

$ws.Range("D39").Value = "'7.24"
$ws.Range("E39").Value = "  -8.12%  "

$ws.Range("E40").Value = "  -9.58%  "

$ws.Range("E41").Value = "  -7.11%  "

$ws.Range("D42").Value = "'0.869"
$ws.Range("E42").Value = "  -4.93%  "

$ws.Range("D43").Value = "'23.97"
$ws.Range("E43").Value = "  -1.64%  "

$ws.Range("E44").Value = "  -1.65%  "

$ws.Range("E45").Value = "  -4.76%  "

$ws.Range("E46").Value = "  -7.55%  "

$ws.Range("E47").Value = "  -6.42%  "

$ws.Range("D48").Value = "'51.97"
$ws.Range("E48").Value = "  -3.84%  "

$ws.Range("D49").Value = "'2.07"
$ws.Range("E49").Value = "  -6.59%  "

$ws.Range("E50").Value = "  -7.28%  "

$ws.Range("D51").Value = "'3.01"
$ws.Range("E51").Value = "  -10.79%  "

